$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.192.96"
$ws.Range("E2").Value = "'  +1.19%  "
$ws.Range("D3").Value = "'2.356.27"
$ws.Range("E3").Value = "'  +2.53%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'302.78"
$ws.Range("E5").Value = "'  +0.54%  "
$ws.Range("D6").Value = "'95.62"
$ws.Range("E6").Value = "'  -0.35%  "
$ws.Range("D7").Value = "'0.505"
$ws.Range("E7").Value = "'  -1.29%  "
$ws.Range("E8").Value = "'  +0.02%  "
$ws.Range("D9").Value = "'0.498"
$ws.Range("E9").Value = "'  +0.80%  "
$ws.Range("D10").Value = "'34.15"
$ws.Range("E10").Value = "'  -1.78%  "
$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "'  +0.12%  "
$ws.Range("D12").Value = "'18.67"
$ws.Range("E12").Value = "'  -2.80%  "
$ws.Range("E13").Value = "'  +3.17%  "
$ws.Range("E14").Value = "'  -0.70%  "
$ws.Range("D15").Value = "'2.723.34"
$ws.Range("E15").Value = "'  +2.72%  "
$ws.Range("D16").Value = "'2.375.16"
$ws.Range("E16").Value = "'  +2.48%  "
$ws.Range("D17").Value = "'0.798"
$ws.Range("E17").Value = "'  +1.76%  "
$ws.Range("D18").Value = "'43.176.96"
$ws.Range("E18").Value = "'  +1.32%  "
$ws.Range("D19").Value = "'12.19"
$ws.Range("E19").Value = "'  -0.98%  "
$ws.Range("E20").Value = "'  +3.88%  "
$ws.Range("E21").Value = "'  -0.01%  "
$ws.Range("E22").Value = "'  +0.57%  "
$ws.Range("D23").Value = "'235.55"
$ws.Range("E23").Value = "'  +0.12%  "
$ws.Range("E24").Value = "'  -1.68%  "
$ws.Range("E25").Value = "'  -0.14%  "
$ws.Range("D26").Value = "'2.42"
$ws.Range("E26").Value = "'  +0.86%  "
$ws.Range("D27").Value = "'24.57"
$ws.Range("E27").Value = "'  -0.20%  "
$ws.Range("E28").Value = "'  +14.79%  "
$ws.Range("D29").Value = "'9.15"
$ws.Range("E29").Value = "'  +0.89%  "
$ws.Range("D30").Value = "'31.33"
$ws.Range("E30").Value = "'  -2.60%  "
$ws.Range("E32").Value = "'  +0.99%  "
$ws.Range("D33").Value = "'0.0725"
$ws.Range("E33").Value = "'  +3.23%  "
$ws.Range("D34").Value = "'17.20"
$ws.Range("E34").Value = "'  -1.94%  "
$ws.Range("D35").Value = "'1.84"
$ws.Range("E35").Value = "'  +5.14%  "
$ws.Range("D36").Value = "'4.38"
$ws.Range("E36").Value = "'  -1.50%  "
$ws.Range("D37").Value = "'2.31"
$ws.Range("E37").Value = "'  -0.73%  "
$ws.Range("D38").Value = "'0.100"
$ws.Range("E38").Value = "'  +0.13%  "
$ws.Range("D39").Value = "'22.46"
$ws.Range("E39").Value = "'  +13.36%  "
$ws.Range("D40").Value = "'2.76"
$ws.Range("E40").Value = "'  +1.51%  "
$ws.Range("E41").Value = "'  -0.33%  "
$ws.Range("D42").Value = "'110.72"
$ws.Range("E42").Value = "'  -32.83%  "
$ws.Range("D43").Value = "'1.942.53"
$ws.Range("E43").Value = "'  -1.33%  "
$ws.Range("D44").Value = "'0.0280"
$ws.Range("E44").Value = "'  +0.64%  "
$ws.Range("E45").Value = "'  +3.50%  "
$ws.Range("D46").Value = "'9.41"
$ws.Range("E46").Value = "'  -10.00%  "
$ws.Range("D47").Value = "'2.75"
$ws.Range("E47").Value = "'  -0.41%  "
$ws.Range("D48").Value = "'2.587.09"
$ws.Range("E48").Value = "'  +2.53%  "
$ws.Range("D49").Value = "'52.89"
$ws.Range("E49").Value = "'  -0.35%  "
$ws.Range("E50").Value = "'  -3.07%  "
$ws.Range("D51").Value = "'72.13"
$ws.Range("E51").Value = "'  +0.79%  "
